$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value = 201
$ws1.Range("F6").Value = 342
$ws1.Range("F7").Value = 225
$ws1.Range("F8").Value = 2179
$ws1.Range("F10").Value = 5356
$ws1.Range("F11").Value = 121

# Sheet "全部类型" (sheet4)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F6").Value = 201
$ws4.Range("F7").Value = 342
$ws4.Range("F8").Value = 225
$ws4.Range("F11").Value = 2179
$ws4.Range("F13").Value = 5356
$ws4.Range("F14").Value = 121
